# Weekly update: insert a new price-report row for "Sandia" (Vega Monumental
# Concepción) above the existing historical rows, pushing rows 152:199 down
# to 153:200 (dimension grows from A1:R199 to A1:R200).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152; Excel shifts rows 152:199 -> 153:200
# and copies formatting (incl. the date style on column D) down with them.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with this week's entry. Columns that
# are constant across this whole market/category block (A, B, C, E, F, G, H,
# Q, R) are copied from the row directly below (the shifted former row 152),
# while the reported fields (D, I, J, K, L, M, N, O, P) get the new values.
$ws.Cells.Item(152, 1).Value  = $ws.Cells.Item(153, 1).Value()   # Mercado ID
$ws.Cells.Item(152, 2).Value  = $ws.Cells.Item(153, 2).Value()   # Mercado
$ws.Cells.Item(152, 3).Value  = $ws.Cells.Item(153, 3).Value()   # Región
$ws.Cells.Item(152, 4).Value  = 45204                            # Fecha
$ws.Cells.Item(152, 5).Value  = $ws.Cells.Item(153, 5).Value()   # Codreg
$ws.Cells.Item(152, 6).Value  = $ws.Cells.Item(153, 6).Value()   # Categoría ID
$ws.Cells.Item(152, 7).Value  = $ws.Cells.Item(153, 7).Value()   # Categoría
$ws.Cells.Item(152, 8).Value  = $ws.Cells.Item(153, 8).Value()   # Variedad
$ws.Cells.Item(152, 9).Value  = "Primera"                        # Calidad
$ws.Cells.Item(152, 10).Value = 500                               # Volumen
$ws.Cells.Item(152, 11).Value = 700                               # Precio mínimo
$ws.Cells.Item(152, 12).Value = 750                               # Precio máximo
$ws.Cells.Item(152, 13).Value = 730                               # Precio promedio ponderado
$ws.Cells.Item(152, 14).Value = "`$/kilo (volumen en unidades)"   # Unidad de comercialización
$ws.Cells.Item(152, 15).Value = "Perú"                            # Origen
$ws.Cells.Item(152, 16).Value = 730                               # Precio $/Kg
$ws.Cells.Item(152, 17).Value = $ws.Cells.Item(153, 17).Value()  # Kg o Unidades
$ws.Cells.Item(152, 18).Value = $ws.Cells.Item(153, 18).Value()  # Clasificación
